$d = $word.ActiveDocument

$pairs = @(
    @("94×11=", "48×29="),
    @("26×94=", "19×68="),
    @("59×63=", "12×51="),
    @("22×49=", "60×98="),
    @("91×71=", "27×66="),
    @("42×54=", "46×95="),
    @("36×31=", "85×13="),
    @("78×44=", "45×96="),
    @("54×62=", "47×35="),
    @("46×14=", "14×90="),
    @("75×67=", "68×64="),
    @("32×31=", "26×83="),
    @("93×93=", "41×65="),
    @("69×99=", "52×99="),
    @("84×71=", "34×85="),
    @("49×64=", "98×37="),
    @("27×55=", "46×54="),
    @("13×99=", "64×34="),
    @("85×45=", "20×32="),
    @("54×53=", "45×62="),
    @("25×75=", "14×39="),
    @("98×52=", "81×54="),
    @("30×82=", "46×15="),
    @("25×50=", "88×47="),
    @("96×19=", "78×88=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
